# Update "想去人数" (want-to-go count) figures on the two sheets that hold
# the full event listing: "展览" (sheet1) and "全部类型" (sheet4).
# F2: 788 -> 792
# F3: 60  -> 61

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 792
    $ws.Range("F3").Value = 61
}
